# Fruta / hortaliza, semanal
# Rewrites the weekly Espárragos price rows (2-8, 11-15) with the updated
# values from the source feed. Columns A,B,C,E,F,G,H,R are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $Fecha, $Calidad, $Volumen, $PrecioMin, $PrecioMax, $PrecioProm, $Unidad, $Origen, $PrecioKg, $KgUnidades) {
    $ws.Cells.Item($Row, 4).Value  = $Fecha        # D - Fecha
    $ws.Cells.Item($Row, 9).Value  = $Calidad       # I - Calidad
    $ws.Cells.Item($Row, 10).Value = $Volumen       # J - Volumen
    $ws.Cells.Item($Row, 11).Value = $PrecioMin     # K - Precio minimo
    $ws.Cells.Item($Row, 12).Value = $PrecioMax     # L - Precio maximo
    $ws.Cells.Item($Row, 13).Value = $PrecioProm    # M - Precio promedio ponderado
    $ws.Cells.Item($Row, 14).Value = $Unidad        # N - Unidad de comercializacion
    $ws.Cells.Item($Row, 15).Value = $Origen        # O - Origen
    $ws.Cells.Item($Row, 16).Value = $PrecioKg      # P - Precio $/Kg
    $ws.Cells.Item($Row, 17).Value = $KgUnidades    # Q - Kg o Unidades
}

Set-Row 2  44160 "Banquete" 210  13000 13000 13000 "$/bandeja 10 kilos" "Región Metropolitana" 1300 10
Set-Row 3  44160 "Primera"  340  11000 11000 11000 "$/bandeja 10 kilos" "Región Metropolitana" 1100 10
Set-Row 4  44160 "Primera"  4300 1200  1200  1200  "$/kilo"              "Región Metropolitana" 1200 1
Set-Row 5  44160 "Segunda"  250  9000  9000  9000  "$/bandeja 10 kilos" "Región Metropolitana" 900  10
Set-Row 6  44160 "Segunda"  2500 1000  1000  1000  "$/kilo"              "Región Metropolitana" 1000 1

Set-Row 7  44162 "Primera"  5200 1100  1100  1100  "$/kilo"              "Provincia de Linares"  1100 1
Set-Row 8  44162 "Segunda"  3400 900   900   900   "$/kilo"              "Provincia de Linares"  900  1

Set-Row 11 44161 "Primera"  4300 1000  1000  1000  "$/kilo"              "Provincia de Linares"  1000 1
Set-Row 12 44161 "Segunda"  2500 800   800   800   "$/kilo"              "Provincia de Linares"  800  1
Set-Row 13 44167 "Primera"  250  13000 13000 13000 "$/bandeja 10 kilos" "Provincia de Linares"  1300 10
Set-Row 14 44167 "Primera"  160  11000 11000 11000 "$/caja 10 kilos"    "Provincia de Linares"  1100 10
Set-Row 15 44167 "Segunda"  120  10000 10000 10000 "$/bandeja 10 kilos" "Provincia de Linares"  1000 10
